$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '''60.684.44'
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = '''  +5.72%  '
$ws.Cells.Item(2, 5).Style = "Normal"

$ws.Cells.Item(3, 4).Value = '''2.637.94'
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = '''  +8.75%  '
$ws.Cells.Item(3, 5).Style = "Normal"

$ws.Cells.Item(4, 4).Value = '''0.997'
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = '''  -0.30%  '
$ws.Cells.Item(4, 5).Style = "Normal"

$ws.Cells.Item(5, 4).Value = '''507.73'
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '''  +3.81%  '
$ws.Cells.Item(5, 5).Style = "Normal"

$ws.Cells.Item(6, 4).Value = '''156.89'
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = '''  +1.71%  '
$ws.Cells.Item(6, 5).Style = "Normal"

$ws.Cells.Item(7, 4).Value = '''0.996'
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = '''  -0.03%  '
$ws.Cells.Item(7, 5).Style = "Normal"

$ws.Cells.Item(8, 4).Value = '''0.588'
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = '''  -4.40%  '
$ws.Cells.Item(8, 5).Style = "Normal"

$ws.Cells.Item(9, 4).Value = '''2.634.07'
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = '''  +7.63%  '
$ws.Cells.Item(9, 5).Style = "Normal"

$ws.Cells.Item(10, 4).Value = '''6.40'
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = '''  +3.49%  '
$ws.Cells.Item(10, 5).Style = "Normal"

$ws.Cells.Item(11, 4).Value = '''0.105'
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = '''  +4.49%  '
$ws.Cells.Item(11, 5).Style = "Normal"

$ws.Cells.Item(12, 5).Value = '''  +2.96%  '
$ws.Cells.Item(12, 5).Style = "Normal"

$ws.Cells.Item(13, 5).Value = '''  +1.12%  '
$ws.Cells.Item(13, 5).Style = "Normal"

$ws.Cells.Item(14, 4).Value = '''3.064.90'
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = '''  +7.38%  '
$ws.Cells.Item(14, 5).Style = "Normal"

$ws.Cells.Item(15, 4).Value = '''60.697.41'
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = '''  +5.91%  '
$ws.Cells.Item(15, 5).Style = "Normal"

$ws.Cells.Item(16, 4).Value = '''21.74'
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = '''  +5.60%  '
$ws.Cells.Item(16, 5).Style = "Normal"

$ws.Cells.Item(17, 5).Value = '''  +5.18%  '
$ws.Cells.Item(17, 5).Style = "Normal"

$ws.Cells.Item(18, 4).Value = '''2.623.52'
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = '''  +7.22%  '
$ws.Cells.Item(18, 5).Style = "Normal"

$ws.Cells.Item(19, 5).Value = '''  +3.49%  '
$ws.Cells.Item(19, 5).Style = "Normal"

$ws.Cells.Item(20, 4).Value = '''344.44'
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = '''  +6.25%  '
$ws.Cells.Item(20, 5).Style = "Normal"

$ws.Cells.Item(21, 4).Value = '''10.43'
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = '''  +4.17%  '
$ws.Cells.Item(21, 5).Style = "Normal"

$ws.Cells.Item(22, 4).Value = '''6.17'
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = '''  +4.30%  '
$ws.Cells.Item(22, 5).Style = "Normal"

$ws.Cells.Item(23, 4).Value = '''1.00'
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = '''  +0.06%  '
$ws.Cells.Item(23, 5).Style = "Normal"

$ws.Cells.Item(24, 4).Value = '''5.77'
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = '''  +0.16%  '
$ws.Cells.Item(24, 5).Style = "Normal"

$ws.Cells.Item(25, 4).Value = '''60.47'
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = '''  +4.52%  '
$ws.Cells.Item(25, 5).Style = "Normal"

$ws.Cells.Item(26, 5).Value = '''  +5.72%  '
$ws.Cells.Item(26, 5).Style = "Normal"

$ws.Cells.Item(27, 5).Value = '''  +3.58%  '
$ws.Cells.Item(27, 5).Style = "Normal"

$ws.Cells.Item(28, 5).Value = '''  -0.47%  '
$ws.Cells.Item(28, 5).Style = "Normal"

$ws.Cells.Item(29, 5).Value = '''  +8.95%  '
$ws.Cells.Item(29, 5).Style = "Normal"

$ws.Cells.Item(30, 4).Value = '''7.56'
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = '''  +4.06%  '
$ws.Cells.Item(30, 5).Style = "Normal"

$ws.Cells.Item(31, 5).Value = '''  -0.22%  '
$ws.Cells.Item(31, 5).Style = "Normal"

$ws.Cells.Item(32, 4).Value = '''156.33'
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = '''  +4.10%  '
$ws.Cells.Item(32, 5).Style = "Normal"

$ws.Cells.Item(33, 5).Value = '''  +3.62%  '
$ws.Cells.Item(33, 5).Style = "Normal"

$ws.Cells.Item(34, 5).Value = '''  +3.37%  '
$ws.Cells.Item(34, 5).Style = "Normal"

$ws.Cells.Item(35, 5).Value = '''  +7.69%  '
$ws.Cells.Item(35, 5).Style = "Normal"

$ws.Cells.Item(36, 5).Value = '''  +6.04%  '
$ws.Cells.Item(36, 5).Style = "Normal"

$ws.Cells.Item(37, 5).Value = '''  +6.50%  '
$ws.Cells.Item(37, 5).Style = "Normal"

$ws.Cells.Item(38, 4).Value = '''309.17'
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = '''  +8.23%  '
$ws.Cells.Item(38, 5).Style = "Normal"

$ws.Cells.Item(39, 5).Value = '''  +4.30%  '
$ws.Cells.Item(39, 5).Style = "Normal"

$ws.Cells.Item(40, 4).Value = '''3.79'
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = '''  +7.57%  '
$ws.Cells.Item(40, 5).Style = "Normal"

$ws.Cells.Item(41, 4).Value = '''1.47'
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = '''  +7.27%  '
$ws.Cells.Item(41, 5).Style = "Normal"

$ws.Cells.Item(42, 5).Value = '''  +27.52%  '
$ws.Cells.Item(42, 5).Style = "Normal"

$ws.Cells.Item(43, 4).Value = '''35.69'
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = '''  +4.85%  '
$ws.Cells.Item(43, 5).Style = "Normal"

$ws.Cells.Item(44, 5).Value = '''  +3.89%  '
$ws.Cells.Item(44, 5).Style = "Normal"

$ws.Cells.Item(45, 4).Value = '''0.0571'
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = '''  +7.36%  '
$ws.Cells.Item(45, 5).Style = "Normal"

$ws.Cells.Item(46, 5).Value = '''  -0.33%  '
$ws.Cells.Item(46, 5).Style = "Normal"

$ws.Cells.Item(47, 4).Value = '''0.999'
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = '''  +0.17%  '
$ws.Cells.Item(47, 5).Style = "Normal"

$ws.Cells.Item(48, 4).Value = '''19.86'
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = '''  +12.98%  '
$ws.Cells.Item(48, 5).Style = "Normal"

$ws.Cells.Item(49, 4).Value = '''4.87'
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = '''  +7.15%  '
$ws.Cells.Item(49, 5).Style = "Normal"

$ws.Cells.Item(50, 5).Value = '''  +3.98%  '
$ws.Cells.Item(50, 5).Style = "Normal"

$ws.Cells.Item(51, 4).Value = '''2.046.58'
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = '''  +7.51%  '
$ws.Cells.Item(51, 5).Style = "Normal"
